$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Diameter of binary tree"
$ws.Range("H25").Value = "DiameterOfBT"

$ws.Range("H25").Select()
